# snc-connector-test-data.xlsx
# "update_time of snc-test of clickhouse was string type update to use"
#
# The clickhouse update_time condition used to be tested with several
# string-formatted date variants (rows 76-79). Those extra variants are
# removed, and the remaining clickhouse update_time test now uses a
# (long/unix-timestamp) value instead of a date string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76 (snc-connector-test-clickhouse-time-1): switch the condition value
# from a quoted date-time string to an epoch/unix timestamp.
$ws.Range("C76").Value = "update_time<'1656038990'"

# Rows 77-79 (snc-connector-test-clickhouse-time-2/3/4) are no longer
# needed - delete them, shifting everything below up by 3 rows.
$ws.Rows("77:79").Delete()

# Put the selection where Excel would naturally leave it after this edit.
$ws.Range("C76").Select()
